# "ya coincide el id de tarjeta de perosnal con la RN3"
#
# The department catalog's "iddepartamento" (col A) and "idsuperior" (col C)
# values are renumbered so the department id matches the staff-card id
# required by business rule RN3. Because the cells already carry a Text
# number format (style index 1 / numFmtId 49), assigning the new ids as
# strings keeps them stored as text (shared strings) rather than numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> (new iddepartamento, new idsuperior)
$updates = @(
    @{ Row = 2;  Id = "10"; Superior = "10" },
    @{ Row = 3;  Id = "11"; Superior = "10" },
    @{ Row = 4;  Id = "2";  Superior = "11" },
    @{ Row = 5;  Id = "7";  Superior = "11" },
    @{ Row = 6;  Id = "8";  Superior = "11" },
    @{ Row = 7;  Id = "5";  Superior = "11" },
    @{ Row = 8;  Id = "12"; Superior = "11" },
    @{ Row = 9;  Id = "13"; Superior = "11" },
    @{ Row = 10; Id = "14"; Superior = "11" },
    @{ Row = 11; Id = "6";  Superior = "10" },
    @{ Row = 12; Id = "15"; Superior = "6"  },
    @{ Row = 13; Id = "16"; Superior = "6"  },
    @{ Row = 14; Id = "17"; Superior = "10" },
    @{ Row = 15; Id = "18"; Superior = "17" },
    @{ Row = 16; Id = "19"; Superior = "17" },
    @{ Row = 17; Id = "20"; Superior = "17" },
    @{ Row = 18; Id = "21"; Superior = "17" },
    @{ Row = 19; Id = "22"; Superior = "10" },
    @{ Row = 20; Id = "23"; Superior = "22" },
    @{ Row = 21; Id = "24"; Superior = "22" },
    @{ Row = 22; Id = "25"; Superior = "22" },
    @{ Row = 23; Id = "26"; Superior = "10" },
    @{ Row = 24; Id = "3";  Superior = "10" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 1).Value = $u.Id        # column A: iddepartamento
    $ws.Cells.Item($u.Row, 3).Value = $u.Superior   # column C: idsuperior
}

# The renumbered ids are now up to two characters wide; re-fit the visible
# columns like Excel does after the values change.
$ws.Columns("A:D").AutoFit() | Out-Null

# Leave the selection where the author ended up after the edit.
$ws.Range("C23").Select() | Out-Null
